# Update cryptocurrency price/volume figures to the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.676.51"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.585.57"
$ws.Range("E3").Value = "  -3.09%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'206.70"
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("E6").Value = "  -3.24%  "
$ws.Range("D8").Value = "'22.29"
$ws.Range("E8").Value = "  -4.75%  "
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("D11").Value = "'0.0869"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("E12").Value = "  -3.13%  "
$ws.Range("D13").Value = "1.580.19"
$ws.Range("E13").Value = "  -4.10%  "
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("D15").Value = "'0.532"
$ws.Range("E15").Value = "  -5.59%  "
$ws.Range("D16").Value = "27.661.26"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("E17").Value = "  -3.59%  "
$ws.Range("D18").Value = "'218.73"
$ws.Range("E18").Value = "  -4.39%  "
$ws.Range("D19").Value = "0.0₃0693"
$ws.Range("E19").Value = "  -3.64%  "
$ws.Range("E20").Value = "  -4.96%  "
$ws.Range("E22").Value = "  -4.99%  "
$ws.Range("D23").Value = "'9.54"
$ws.Range("E23").Value = "  -5.35%  "
$ws.Range("D24").Value = "'1.97"
$ws.Range("E24").Value = "  -5.08%  "
$ws.Range("D25").Value = "'153.56"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D28").Value = "'15.09"
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("E29").Value = "  -4.00%  "
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("E31").Value = "  -3.38%  "
$ws.Range("E32").Value = "  -5.51%  "
$ws.Range("D33").Value = "1.380.94"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("E34").Value = "  -5.00%  "
$ws.Range("E35").Value = "  -5.44%  "
$ws.Range("D36").Value = "'0.968"
$ws.Range("E36").Value = "  -4.28%  "
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("D39").Value = "'0.541"
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("D43").Value = "'1.79"
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("D44").Value = "'63.88"
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("E46").Value = "  -4.24%  "
$ws.Range("D47").Value = "1.720.98"
$ws.Range("E47").Value = "  -3.14%  "
$ws.Range("D48").Value = "'87.70"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("E50").Value = "  -5.25%  "
$ws.Range("D51").Value = "'0.0499"
$ws.Range("E51").Value = "  -1.18%  "
